$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1744
$ws.Range("I40").Value = 1714.2858
$ws.Range("J40").Value = 1778.6666
$ws.Range("K40").Value = 1714.2858
$ws.Range("L40").Value = 1778.6666
$ws.Range("M40").Value = -1539.2858
$ws.Range("N40").Value = -2128.6666

$ws.Range("H64").Value = 3356.3845
$ws.Range("I64").Value = 3099.875
$ws.Range("J64").Value = 3422.5806
$ws.Range("K64").Value = 3099.875
$ws.Range("L64").Value = 3422.5806
$ws.Range("M64").Value = -2851.875
$ws.Range("N64").Value = -3918.5806

$ws.Range("H67").Value = 3356.3845
$ws.Range("I67").Value = 3099.875
$ws.Range("J67").Value = 3422.5806
$ws.Range("K67").Value = 3099.875
$ws.Range("L67").Value = 3422.5806
$ws.Range("M67").Value = -2241.875
$ws.Range("N67").Value = -5138.580599999999

$ws.Range("H74").Value = 4700.6
$ws.Range("I74").Value = 4251.5
$ws.Range("K74").Value = 4251.5
$ws.Range("M74").Value = -3315.5

$ws.Range("H77").Value = 4700.6
$ws.Range("I77").Value = 4251.5
$ws.Range("K77").Value = 21257.5
$ws.Range("M77").Value = -16577.5

$ws.Range("H132").Value = 1791.6613
$ws.Range("I132").Value = 1527.7174
$ws.Range("J132").Value = 2550.5
$ws.Range("K132").Value = 4583.1522
$ws.Range("L132").Value = 7651.5
$ws.Range("M132").Value = -2053.1522
$ws.Range("N132").Value = -12711.5

$ws.Range("H138").Value = 3123.3542
$ws.Range("I138").Value = 1968.2759
$ws.Range("J138").Value = 4886.3687
$ws.Range("K138").Value = 5904.8277
$ws.Range("L138").Value = 14659.1061
$ws.Range("M138").Value = -764.8276999999998
$ws.Range("N138").Value = -24939.1061

$ws.Range("H141").Value = 4265.725
$ws.Range("I141").Value = 1698.0278
$ws.Range("J141").Value = 27375
$ws.Range("K141").Value = 5094.0834
$ws.Range("L141").Value = 82125
$ws.Range("M141").Value = 85.91659999999956
$ws.Range("N141").Value = -92485

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8230.522999999999
$ws.Range("I32").Value = 8086.4287
$ws.Range("J32").Value = 9815.571
$ws.Range("K32").Value = 8086.4287
$ws.Range("L32").Value = 9815.571
$ws.Range("M32").Value = -7799.4287
$ws.Range("N32").Value = -10389.571

$ws.Range("H63").Value = 4417.3335
$ws.Range("I63").Value = 4626.25
$ws.Range("K63").Value = 4626.25
$ws.Range("M63").Value = -3940.25

$ws.Range("H66").Value = 4417.3335
$ws.Range("I66").Value = 4626.25
$ws.Range("K66").Value = 23131.25
$ws.Range("M66").Value = -19699.25

$ws.Range("H88").Value = 3081.6667
$ws.Range("J88").Value = 3118
$ws.Range("L88").Value = 3118
$ws.Range("N88").Value = -3930

$ws.Range("H91").Value = 3081.6667
$ws.Range("J91").Value = 3118
$ws.Range("L91").Value = 3118
$ws.Range("N91").Value = -5926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 193333.33
$ws.Range("I82").Value = 193333.33
$ws.Range("K82").Value = 193333.33
$ws.Range("M82").Value = -192950.33

$ws.Range("H85").Value = 193333.33
$ws.Range("I85").Value = 193333.33
$ws.Range("K85").Value = 193333.33
$ws.Range("M85").Value = -192007.33

$ws.Range("H86").Value = 145937
$ws.Range("I86").Value = 4258.7144
$ws.Range("J86").Value = 287615.28
$ws.Range("K86").Value = 4258.7144
$ws.Range("L86").Value = 287615.28
$ws.Range("M86").Value = -3135.7144
$ws.Range("N86").Value = -289861.28

$ws.Range("H89").Value = 145937
$ws.Range("I89").Value = 4258.7144
$ws.Range("J89").Value = 287615.28
$ws.Range("K89").Value = 21293.572
$ws.Range("L89").Value = 1438076.4
$ws.Range("M89").Value = -15677.572
$ws.Range("N89").Value = -1449308.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 226399.84
$ws.Range("I132").Value = 301186.34
$ws.Range("J132").Value = 2040.4
$ws.Range("K132").Value = 903559.02
$ws.Range("L132").Value = 6121.200000000001
$ws.Range("M132").Value = -901029.02
$ws.Range("N132").Value = -11181.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 3286.0625
$ws.Range("I112").Value = 849.25
$ws.Range("J112").Value = 4098.3335
$ws.Range("K112").Value = 2547.75
$ws.Range("L112").Value = 12295.0005
$ws.Range("M112").Value = -1439.75
$ws.Range("N112").Value = -14511.0005

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3316.1538
$ws.Range("I80").Value = 3563.75
$ws.Range("J80").Value = 2920
$ws.Range("K80").Value = 3563.75
$ws.Range("L80").Value = 2920
$ws.Range("M80").Value = -2565.75
$ws.Range("N80").Value = -4916

$ws.Range("H83").Value = 3316.1538
$ws.Range("I83").Value = 3563.75
$ws.Range("J83").Value = 2920
$ws.Range("K83").Value = 17818.75
$ws.Range("L83").Value = 14600
$ws.Range("M83").Value = -12826.75
$ws.Range("N83").Value = -24584

$ws.Range("H131").Value = 36481.5
$ws.Range("J131").Value = 36481.5
$ws.Range("L131").Value = 36481.5
$ws.Range("N131").Value = -46561.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 828.6667

$ws.Range("H27").Value = 828.6667

$ws.Range("H46").Value = 916.6667
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 1750
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 1750
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -2126

$ws.Range("H68").Value = 3440.2
$ws.Range("I68").Value = 2833.1667
$ws.Range("J68").Value = 4350.75
$ws.Range("K68").Value = 2833.1667
$ws.Range("L68").Value = 4350.75
$ws.Range("M68").Value = -2084.1667
$ws.Range("N68").Value = -5848.75

$ws.Range("H71").Value = 3440.2
$ws.Range("I71").Value = 2833.1667
$ws.Range("J71").Value = 4350.75
$ws.Range("K71").Value = 14165.8335
$ws.Range("L71").Value = 21753.75
$ws.Range("M71").Value = -10421.8335
$ws.Range("N71").Value = -29241.75

$ws.Range("H82").Value = 2422.75
$ws.Range("I82").Value = 1676.1
$ws.Range("K82").Value = 1676.1
$ws.Range("M82").Value = -1315.1

$ws.Range("H85").Value = 2422.75
$ws.Range("I85").Value = 1676.1
$ws.Range("K85").Value = 1676.1
$ws.Range("M85").Value = -428.0999999999999

$ws.Range("H132").Value = 4076.4583
$ws.Range("I132").Value = 3274
$ws.Range("J132").Value = 6483.8335
$ws.Range("K132").Value = 9822
$ws.Range("L132").Value = 19451.5005
$ws.Range("M132").Value = -7292
$ws.Range("N132").Value = -24511.5005

$ws.Range("H136").Value = 3323.75
$ws.Range("I136").Value = 2453
$ws.Range("K136").Value = 7359
$ws.Range("M136").Value = -4809

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J62").Value = 4777.778
$ws.Range("L62").Value = 4777.778
$ws.Range("N62").Value = -6025.778

$ws.Range("J65").Value = 4777.778
$ws.Range("L65").Value = 23888.89
$ws.Range("N65").Value = -30128.89

$ws.Range("H81").Value = 44983.867
$ws.Range("I81").Value = 40181.383
$ws.Range("J81").Value = 76200
$ws.Range("K81").Value = 80362.766
$ws.Range("L81").Value = 152400
$ws.Range("M81").Value = -79301.766
$ws.Range("N81").Value = -154522

$ws.Range("H84").Value = 44983.867
$ws.Range("I84").Value = 40181.383
$ws.Range("J84").Value = 76200
$ws.Range("K84").Value = 401813.83
$ws.Range("L84").Value = 762000
$ws.Range("M84").Value = -396509.83
$ws.Range("N84").Value = -772608

$ws.Range("H123").Value = 24994.428
$ws.Range("J123").Value = 24994.428
$ws.Range("L123").Value = 24994.428
$ws.Range("N123").Value = -34794.428

$ws.Range("H132").Value = 940.2969000000001
$ws.Range("I132").Value = 738.44446
$ws.Range("K132").Value = 2215.33338
$ws.Range("M132").Value = 314.66662
